$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New "Trening" header in F1 -- reuse the bold/bordered/centered style that
#    the other headers (A1:E1) already use, by copy/pasting formats only.
# ---------------------------------------------------------------------------
$ws.Cells.Item(1, 6).Value = "Trening"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Rewrite the data rows. Timestamps move from text strings to real Excel
#    date/time serial values, the data set grows from 6 to 12 rows, and a
#    "Trening" (session part) label is added to every row.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = 45685.64609351852
$ws.Cells.Item(2, 2).Value = 937.4
$ws.Cells.Item(2, 3).Value = 13.98
$ws.Cells.Item(2, 4).Value = 4.269063949584963
$ws.Cells.Item(2, 5).Value = "10-15"
$ws.Cells.Item(2, 6).Value = "Duża Gra"

$ws.Cells.Item(3, 1).Value = 45685.64963055556
$ws.Cells.Item(3, 2).Value = 1243
$ws.Cells.Item(3, 3).Value = 14.32
$ws.Cells.Item(3, 4).Value = 4.132714033126831
$ws.Cells.Item(3, 5).Value = "10-15"
$ws.Cells.Item(3, 6).Value = "Duża Gra"

$ws.Cells.Item(4, 1).Value = 45685.66339328704
$ws.Cells.Item(4, 2).Value = 2432.1
$ws.Cells.Item(4, 3).Value = 14.85
$ws.Cells.Item(4, 4).Value = 4.18912148475647
$ws.Cells.Item(4, 5).Value = "10-15"
$ws.Cells.Item(4, 6).Value = "Duża Gra"

$ws.Cells.Item(5, 1).Value = 45685.6460900463
$ws.Cells.Item(5, 2).Value = 937.1
$ws.Cells.Item(5, 3).Value = 8.74
$ws.Cells.Item(5, 4).Value = 3.658652578081405
$ws.Cells.Item(5, 5).Value = "5-10"
$ws.Cells.Item(5, 6).Value = "Duża Gra"

$ws.Cells.Item(6, 1).Value = 45685.64962708333
$ws.Cells.Item(6, 2).Value = 1242.7
$ws.Cells.Item(6, 3).Value = 9.609999999999999
$ws.Cells.Item(6, 4).Value = 3.690100363322667
$ws.Cells.Item(6, 5).Value = "5-10"
$ws.Cells.Item(6, 6).Value = "Duża Gra"

$ws.Cells.Item(7, 1).Value = 45685.65109236111
$ws.Cells.Item(7, 2).Value = 1369.3
$ws.Cells.Item(7, 3).Value = 9.119999999999999
$ws.Cells.Item(7, 4).Value = 3.482657057898386
$ws.Cells.Item(7, 5).Value = "5-10"
$ws.Cells.Item(7, 6).Value = "Duża Gra"

$ws.Cells.Item(8, 1).Value = 45685.66953564815
$ws.Cells.Item(8, 2).Value = 2962.8
$ws.Cells.Item(8, 3).Value = 14.66
$ws.Cells.Item(8, 4).Value = 3.711813654218403
$ws.Cells.Item(8, 5).Value = "10-15"
$ws.Cells.Item(8, 6).Value = "Mała Gra"

$ws.Cells.Item(9, 1).Value = 45685.67883310185
$ws.Cells.Item(9, 2).Value = 3766.1
$ws.Cells.Item(9, 3).Value = 11.49
$ws.Cells.Item(9, 4).Value = 3.325542994907926
$ws.Cells.Item(9, 5).Value = "10-15"
$ws.Cells.Item(9, 6).Value = "Mała Gra"

$ws.Cells.Item(10, 1).Value = 45685.68071273148
$ws.Cells.Item(10, 2).Value = 3928.5
$ws.Cells.Item(10, 3).Value = 12.71
$ws.Cells.Item(10, 4).Value = 3.427133185522896
$ws.Cells.Item(10, 5).Value = "10-15"
$ws.Cells.Item(10, 6).Value = "Mała Gra"

$ws.Cells.Item(11, 1).Value = 45685.67769768518
$ws.Cells.Item(11, 2).Value = 3668
$ws.Cells.Item(11, 3).Value = 9.789999999999999
$ws.Cells.Item(11, 4).Value = 2.70002692086356
$ws.Cells.Item(11, 5).Value = "5-10"
$ws.Cells.Item(11, 6).Value = "Mała Gra"

$ws.Cells.Item(12, 1).Value = 45685.67883078704
$ws.Cells.Item(12, 2).Value = 3765.9
$ws.Cells.Item(12, 3).Value = 8.57
$ws.Cells.Item(12, 4).Value = 3.03781158583505
$ws.Cells.Item(12, 5).Value = "5-10"
$ws.Cells.Item(12, 6).Value = "Mała Gra"

$ws.Cells.Item(13, 1).Value = 45685.68070925926
$ws.Cells.Item(13, 2).Value = 3928.2
$ws.Cells.Item(13, 3).Value = 8.609999999999999
$ws.Cells.Item(13, 4).Value = 3.017870221819195
$ws.Cells.Item(13, 5).Value = "5-10"
$ws.Cells.Item(13, 6).Value = "Mała Gra"

# ---------------------------------------------------------------------------
# 3) Apply the "YYYY-MM-DD HH:MM:SS" date/time number format to column A
#    (rows 2-13). A throwaway cell is used first to register the lower-case
#    format variant exactly as the source workbook's style table records it
#    (numFmtId 164 = "yyyy-mm-dd h:mm:ss", numFmtId 165 = the upper-case one
#    actually applied to the cells), then cleared so it leaves no trace in
#    the sheet data / used range.
# ---------------------------------------------------------------------------
$scratch = $ws.Cells.Item(1, 26)
$scratch.NumberFormat = "yyyy-mm-dd h:mm:ss"
$scratch.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$scratch.Clear()

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
